# Refresh the cryptocurrency Price (D) and Volume(1h) (E) columns with the
# latest scraped figures. Price values are stored as plain text (some use
# "."-grouped formatting like "23.647.39" that is not a valid number), so we
# force a Text number format before writing them -- otherwise Excel
# auto-converts plain decimals (e.g. "0.9980") into numeric values and the
# literal formatting is lost. The format flag is cleared again afterwards so
# the cell keeps its original (default) style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.647.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.647.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9980"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9983"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "304.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3801"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.98"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3606"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.248"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08199"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9978"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.526"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.379"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001231"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.637.82"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "96.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06968"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.734"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9981"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.640.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.520"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.05%  "
$ws.Range("E26").Value = "  -2.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("E28").Value = "  +1.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.210"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.833.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.778"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.087"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.050"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -9.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02806"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2517"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08816"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.098"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.07034"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7064"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.331"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.93"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6510"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.339"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9981"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.978"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07981"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "127.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.190"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.66%  "
